# Apply edit: add rows 26-49 to 'Test Results' worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('www.varoom.com', 'https://www.varoom.com/property/gorgeous-loft-style-house-w-luxury-charm-comfort/HA-3212129747', 'Hybrid', 'Property available in date range', 'The property ''* Gorgeous Loft Style House w/Luxury, Charm, & Comfort *'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/nice-bedroom-next-fells-point/BC-4875257', 'Hybrid', 'Property available in date range', 'The property ''Nice bedroom next fells point'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/1840s-carrollton-inn/EP-4238812', 'Hybrid', 'Property available in date range', 'The property ''1840s Carrollton Inn'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/rollingside-guest-suite/HA-3212133687', 'Hybrid', 'Property available in date range', 'The property ''Rollingside Guest Suite'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/chic-3-bedroom-house-in-baltimore/BC-12836036', 'Hybrid', 'Property available in date range', 'The property ''Chic 3 Bedroom house in Baltimore'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/entire-unit-5min-to-m-tbank-camden-yards/EP-78004082', 'Hybrid', 'Property available in date range', 'The property ''Entire Unit-5min to M&TBank Camden Yards'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/glen-burnie-escape-pet-friendly-with-easy-access-to-multiple-attractions/HA-3213843374', 'Hybrid', 'Property available in date range', 'The property ''Glen Burnie Escape - Pet friendly with easy access to multiple attractions.'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/fells-point-charm-doubled/BC-10671028', 'Hybrid', 'Property available in date range', 'The property ''Fells Point Charm Doubled!'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/room-in-guest-room-twin-bunk-style-bedroom-close-to-downtown/EP-84610916', 'Hybrid', 'Property available in date range', 'The property ''Room in Guest Room - Twin Bunk-style Bedroom Close to Downtown'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/big-3-story-house-near-stadiums-aquarium-museums-5-bedrooms-parking/HA-3213607613', 'Hybrid', 'Property available in date range', 'The property ''Big 3 story house near stadiums, aquarium, museums, 5 bedrooms, parking'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/huge-downtown-baltimore-home-with-parking/BC-13221476', 'Hybrid', 'Property available in date range', 'The property ''Huge Downtown Baltimore home with parking'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/room-in-guest-room-cozy-bedroom-close-to-downtown/EP-84610923', 'Hybrid', 'Property available in date range', 'The property ''Room in Guest Room - Cozy Bedroom Close to Downtown'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/blue-house-in-hampden/HA-3213840308', 'Hybrid', 'Property available in date range', 'The property ''Blue House in Hampden'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/home-away-from-home/BC-12692682', 'Hybrid', 'Property available in date range', 'The property ''Home away from home!'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/guesthouse-by-good-neighbor/EP-96404669', 'Hybrid', 'Property available in date range', 'The property ''guesthouse by good neighbor'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/cozy-1-bedroom-with-ac-and-wifi-in-marvelous-baltimore/HA-3214161410', 'Hybrid', 'Property available in date range', 'The property ''Cozy 1-bedroom with AC and WiFi in marvelous Baltimore.'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/entire-house-king-bed-near-downtown-free-parking/BC-12825912', 'Hybrid', 'Property available in date range', 'The property ''Entire House King Bed near Downtown Free Parking'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/four-seasons-hotel-baltimore/EP-4603198', 'Hybrid', 'Property available in date range', 'The property ''Four Seasons Hotel Baltimore'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/fed-hill-luxe-with-harbor-view-rooftop-deck-parking-ev/HA-3214162843', 'Hybrid', 'Property available in date range', 'The property ''Fed Hill Luxe with Harbor View Rooftop Deck, Parking & EV'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/rachael-s-dowry-bed-and-breakfast/BC-529548', 'Hybrid', 'Property available in date range', 'The property ''Rachael''s Dowry Bed and Breakfast'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/roost-baltimore/EP-95348917', 'Hybrid', 'Property available in date range', 'The property ''ROOST Baltimore'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/tiny-house-24-7-no-contact-self-check-in-bwi-area/HA-3214225097', 'Hybrid', 'Property available in date range', 'The property ''Tiny House 24/7 No Contact Self Check-In BWI area'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/1840s-carrollton-inn/BC-1186462', 'Hybrid', 'Property available in date range', 'The property ''1840s Carrollton Inn'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
    ,@('www.varoom.com', 'https://www.varoom.com/property/hotel-revival-baltimore/EP-28314', 'Hybrid', 'Property available in date range', 'The property ''Hotel Revival Baltimore'' is Available in the specified date range. | **Location:** Cookside | **Date Range:** Check-in: 2025-02-09, Check-out: 2025-02-10')
)

$startRow = 26
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = $row[4]
}
